# Fixed subordinating header for question three.
#
# The answer to "What tools have been the most difficult to learn?" had a
# literal Markdown sub-header ("## Why have they been particularly
# difficult to learn and use?") typed inline inside the answer paragraph
# instead of being a real, subordinate Word heading. This promotes that
# text into its own Heading2 paragraph (matching the "Why?" subordinate
# heading used for question two) and strips the leftover Markdown/text
# from the answer paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that still has the inline "## ... ?" markdown
# header embedded in its text (search by content so this isn't tied to a
# fixed paragraph index).
$markerText = "## Why have they been particularly difficult to learn and use? "

$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains($markerText)) {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find the paragraph containing the inline markdown header"
}

# Insert a brand new, empty paragraph immediately before the answer
# paragraph. Because the insertion point precedes $targetPara's start,
# $targetPara now refers to that new (still empty) paragraph, and the
# original answer text has shifted down into the paragraph right after
# it.
$targetPara.Range.InsertParagraphBefore()

# Fill in the new subordinate heading's text.
$headingInsertionPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$headingInsertionPoint.Text = "Why have they been particularly difficult to learn and use?"

# Give it the same subordinate heading style used elsewhere ("Why?").
$targetPara.Style = "Heading2"

# Stamp it with the same kind of (zero-length, paragraph-leading) bookmark
# Pandoc/Word uses for every other header in this document.
$bookmarkPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$d.Bookmarks.Add("why-have-they-been-particularly-difficult-to-learn-and-use", $bookmarkPoint) | Out-Null

# Remove the now-redundant inline markdown header text from the answer
# paragraph that follows.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute($markerText, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
